$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2128
$ws1.Range("F4").Value = 17
$ws1.Range("F5").Value = 11085
$ws1.Range("F6").Value = 191
$ws1.Range("F10").Value = 10975
$ws1.Range("F14").Value = 1712
$ws1.Range("F15").Value = 5497
$ws1.Range("F17").Value = 3421

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 564

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2128
$ws4.Range("F4").Value = 564
$ws4.Range("F5").Value = 17
$ws4.Range("F7").Value = 11085
$ws4.Range("F8").Value = 191
$ws4.Range("F12").Value = 10975
$ws4.Range("F16").Value = 1712
$ws4.Range("F17").Value = 5497
$ws4.Range("F19").Value = 3421
